$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1 / A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 17:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1648959
$ws.Range("C4").Value = 3865
$ws.Range("E4").Value = 1147915

# Row 14 - India
$ws.Range("B14").Value = 130506
$ws.Range("C14").Value = 5712
$ws.Range("D14").Value = 53947
$ws.Range("E14").Value = 72709
$ws.Range("G14").Value = 124
$ws.Range("H14").Value = 3850

# Rows 19 & 20 - Chile overtakes Mexico in the ranking, so the two
# countries swap rows (row19 becomes Chile, row20 becomes Mexico) and
# Chile's stats are updated while Mexico's stay the same as before.
$ws.Range("A19").Value = "Chile"
$ws.Range("B19").Value = 65393
$ws.Range("C19").Value = 3536
$ws.Range("D19").Value = 26546
$ws.Range("E19").Value = 38174
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 673

$ws.Range("A20").Value = "Mexico"
$ws.Range("B20").Value = 62527
$ws.Range("C20").Value = 2960
$ws.Range("D20").Value = 42725
$ws.Range("E20").Value = 12813
$ws.Range("G20").Value = 479
$ws.Range("H20").Value = 6989

# Row 81 - Tayikistan
$ws.Range("B81").Value = 2738
$ws.Range("C81").Value = 187
$ws.Range("D81").Value = 1223
$ws.Range("E81").Value = 1471
